$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "A4"  = -20.575
    "C4"  = -12.566
    "A6"  = -22.157
    "A7"  = -19.948
    "D7"  = -8.565
    "D8"  = -8.401
    "C9"  = -11
    "D10" = -7.675
    "C12" = -11.109
    "D13" = -7.590999999999999
    "A16" = -22.025
    "D16" = -8.053000000000001
    "C17" = -13.392
    "C18" = -12.457
    "C19" = -12.126
    "A20" = -19.828
    "C20" = -11.729
    "C26" = -12.222
    "A29" = -21.324
    "D30" = -7.139
    "C31" = -13.298
    "A32" = -21.772
    "C39" = -12.133
    "A40" = -19.998
    "C40" = -12.038
    "D40" = -7.996
    "C41" = -12.117
    "C42" = -12.29
    "C43" = -12.216
    "D44" = -7.63
    "A46" = -21.89
    "C47" = -12.082
    "C48" = -11.832
    "A51" = -21.71
    "A52" = -21.972
    "A57" = -22.196
    "A59" = -22.493
    "A62" = -21.97
    "C63" = -11.43
    "C64" = -10.674
    "A66" = -21.672
    "A73" = -20.345
    "A74" = -21.197
    "C76" = -12.644
    "C81" = -12.975
    "C89" = -13.075
    "D89" = -8.128000000000002
    "D91" = -7.580000000000001
    "A92" = -21.759
    "C94" = -11.375
    "A100" = -22.259
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
